$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New labels
$ws.Range("B11").Value = "Process Time "
$ws.Range("A12").Value = "MA"
$ws.Range("A13").Value = "MB"

# Row 12: sum of B2:B4 ... K2:K4 (machine A process times)
$ws.Range("B12").Formula = "=SUM(B2:B4)"
$ws.Range("C12:K12").Formula = "=SUM(C2:C4)"

# Row 13: sum of B5 ... K5 (machine B process times)
$ws.Range("B13").Formula = "=SUM(B5)"
$ws.Range("C13:K13").Formula = "=SUM(C5)"

# Empty formatted cell K15 (underline font applied, no value)
$ws.Range("K15").Font.Underline = $true

# Defined names referring to the new summary table
$wb.Names.Add("precess_time", '=Foglio1!$B$12:$K$13')
$wb.Names.Add("process_time", '=Foglio1!$B$12:$K$13')

# Selection as left by the author
$ws.Range("E35").Select() | Out-Null
